$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns before column C, shifting the existing
# C:N data (Venue .. geometry) right to E:P.
$ws.Range("C1:D1").EntireColumn.Insert()

# Header row: the old "Unnamed: 0" column (B) now has two duplicate
# index columns next to it (C, D) before the shifted "Unnamed: 0" header
# that used to sit in B lands in D.
$ws.Range("B1").Value = "Unnamed: 0.2"
$ws.Range("C1").Value = "Unnamed: 0.1"
$ws.Range("D1").Value = "Unnamed: 0"

# Fill the two new index columns with the same row-index values already
# present in column B for every data row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $idx = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 3).Value = $idx
    $ws.Cells.Item($r, 4).Value = $idx
}
